$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 1645.6364
$ws.Range("I98").Value = 955.7778
$ws.Range("J98").Value = 4750
$ws.Range("K98").Value = 955.7778
$ws.Range("L98").Value = 4750
$ws.Range("M98").Value = 542.2222
$ws.Range("N98").Value = -7746

# Row 122
$ws.Range("H122").Value = 1645.6364
$ws.Range("I122").Value = 955.7778
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 2867.3334
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -417.3334
$ws.Range("N122").Value = -19150

# Row 137
$ws.Range("H137").Value = 2404.0264
$ws.Range("J137").Value = 3393.3157
$ws.Range("L137").Value = 10179.9471
$ws.Range("N137").Value = -15279.9471

# Row 138
$ws.Range("H138").Value = 2947.4
$ws.Range("J138").Value = 3658.48
$ws.Range("L138").Value = 10975.44
$ws.Range("N138").Value = -21255.44

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 42152.4
$ws.Range("I2").Value = 2500
$ws.Range("J2").Value = 52065.5
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 52065.5
$ws.Range("M2").Value = -2387
$ws.Range("N2").Value = -52291.5

# Row 32
$ws.Range("H32").Value = 2614.8723
$ws.Range("I32").Value = 1777.8536
$ws.Range("J32").Value = 8334.5
$ws.Range("K32").Value = 1777.8536
$ws.Range("L32").Value = 8334.5
$ws.Range("M32").Value = -1490.8536
$ws.Range("N32").Value = -8908.5

# Row 45
$ws.Range("H45").Value = 55559090
$ws.Range("I45").Value = 100001100
$ws.Range("J45").Value = 6576.25
$ws.Range("K45").Value = 100001100
$ws.Range("L45").Value = 6576.25
$ws.Range("M45").Value = -100000723
$ws.Range("N45").Value = -7330.25

# Row 61
$ws.Range("H61").Value = 4569
$ws.Range("J61").Value = 7916.3335
$ws.Range("L61").Value = 7916.3335
$ws.Range("N61").Value = -8340.333500000001

# Row 74
$ws.Range("H74").Value = 15152952
$ws.Range("I74").Value = 19609146
$ws.Range("K74").Value = 19609146
$ws.Range("M74").Value = -19608272

# Row 77
$ws.Range("H77").Value = 15152952
$ws.Range("I77").Value = 19609146
$ws.Range("K77").Value = 98045730
$ws.Range("M77").Value = -98041362

# Row 96
$ws.Range("H96").Value = 16714.625
$ws.Range("J96").Value = 16714.625
$ws.Range("L96").Value = 16714.625
$ws.Range("N96").Value = -22206.625

# Row 97
$ws.Range("H97").Value = 741.1667
$ws.Range("I97").Value = 414.8
$ws.Range("K97").Value = 414.8
$ws.Range("M97").Value = 81.19999999999999

# Row 116
$ws.Range("H116").Value = 42152.4
$ws.Range("I116").Value = 2500
$ws.Range("J116").Value = 52065.5
$ws.Range("K116").Value = 2500
$ws.Range("L116").Value = 52065.5
$ws.Range("M116").Value = -206
$ws.Range("N116").Value = -56653.5

# Row 132
$ws.Range("H132").Value = 5157.9165
$ws.Range("I132").Value = 2880.5
$ws.Range("K132").Value = 8641.5
$ws.Range("M132").Value = -6111.5

# Row 136
$ws.Range("H136").Value = 4569
$ws.Range("J136").Value = 7916.3335
$ws.Range("L136").Value = 23749.0005
$ws.Range("N136").Value = -28849.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 42152.4
$ws.Range("I3").Value = 2500
$ws.Range("J3").Value = 52065.5
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 52065.5
$ws.Range("M3").Value = -2386
$ws.Range("N3").Value = -52293.5

# Row 134
$ws.Range("H134").Value = 3001.3142
$ws.Range("I134").Value = 1607.6207
$ws.Range("K134").Value = 4822.8621
$ws.Range("M134").Value = -2287.8621

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 4349.5
$ws.Range("J22").Value = 6800
$ws.Range("L22").Value = 6800
$ws.Range("N22").Value = -7500

# Row 70
$ws.Range("H70").Value = 44545
$ws.Range("J70").Value = 44545
$ws.Range("L70").Value = 44545
$ws.Range("N70").Value = -45175

# Row 73
$ws.Range("H73").Value = 44545
$ws.Range("J73").Value = 44545
$ws.Range("L73").Value = 44545
$ws.Range("N73").Value = -46729

# Row 99
$ws.Range("H99").Value = 2335.4285
$ws.Range("I99").Value = 2369.6
$ws.Range("K99").Value = 2369.6
$ws.Range("M99").Value = -871.5999999999999

# Row 126
$ws.Range("H126").Value = 2335.4285
$ws.Range("I126").Value = 2369.6
$ws.Range("K126").Value = 7108.799999999999
$ws.Range("M126").Value = -4638.799999999999

# Row 134
$ws.Range("H134").Value = 2717.3
$ws.Range("I134").Value = 1950.1305
$ws.Range("K134").Value = 5850.3915
$ws.Range("M134").Value = -3315.3915

$ws = $wb.Worksheets.Item("CUL")
# Row 42
$ws.Range("H42").Value = 11000
$ws.Range("J42").Value = 11000
$ws.Range("L42").Value = 33000
$ws.Range("N42").Value = -34068

# Row 69
$ws.Range("H69").Value = 8502
$ws.Range("J69").Value = 8938.375
$ws.Range("L69").Value = 26815.125
$ws.Range("N69").Value = -28437.125

# Row 72
$ws.Range("H72").Value = 8502
$ws.Range("J72").Value = 8938.375
$ws.Range("L72").Value = 80445.375
$ws.Range("N72").Value = -88557.375

$ws = $wb.Worksheets.Item("GSM")
# Row 55
$ws.Range("H55").Value = 395705.75
$ws.Range("I55").Value = 20000
$ws.Range("K55").Value = 20000
$ws.Range("M55").Value = -19673

# Row 80
$ws.Range("H80").Value = 6423
$ws.Range("J80").Value = 9519.4
$ws.Range("L80").Value = 9519.4
$ws.Range("N80").Value = -11515.4

# Row 83
$ws.Range("H83").Value = 6423
$ws.Range("J83").Value = 9519.4
$ws.Range("L83").Value = 47597
$ws.Range("N83").Value = -57581

# Row 95
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492

# Row 97
$ws.Range("H97").Value = 1598.7826
$ws.Range("I97").Value = 1296.9375
$ws.Range("K97").Value = 1296.9375
$ws.Range("M97").Value = -800.9375

# Row 113
$ws.Range("H113").Value = 3799.7144
$ws.Range("I113").Value = 1866.3334
$ws.Range("J113").Value = 5249.75
$ws.Range("K113").Value = 1866.3334
$ws.Range("L113").Value = 5249.75
$ws.Range("M113").Value = 303.6666
$ws.Range("N113").Value = -9589.75

# Row 122
$ws.Range("H122").Value = 6703.375
$ws.Range("I122").Value = 3724
$ws.Range("K122").Value = 11172
$ws.Range("M122").Value = -8722

# Row 126
$ws.Range("H126").Value = 3309.7
$ws.Range("I126").Value = 1654.7142
$ws.Range("J126").Value = 7171.3335
$ws.Range("K126").Value = 4964.142599999999
$ws.Range("L126").Value = 21514.0005
$ws.Range("M126").Value = -2494.142599999999
$ws.Range("N126").Value = -26454.0005

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# Row 136
$ws.Range("H136").Value = 55375.7
$ws.Range("J136").Value = 55375.7
$ws.Range("L136").Value = 166127.1
$ws.Range("N136").Value = -171227.1

# Row 139
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4230.1875
$ws.Range("I40").Value = 3792.3333
$ws.Range("K40").Value = 3792.3333
$ws.Range("M40").Value = -3656.3333

# Row 122
$ws.Range("H122").Value = 4669.9033
$ws.Range("I122").Value = 4267.846
$ws.Range("J122").Value = 6760.6
$ws.Range("K122").Value = 12803.538
$ws.Range("L122").Value = 20281.8
$ws.Range("M122").Value = -10353.538
$ws.Range("N122").Value = -25181.8

# Row 133
$ws.Range("H133").Value = 57775.668
$ws.Range("J133").Value = 57775.668
$ws.Range("L133").Value = 57775.668
$ws.Range("N133").Value = -62835.668

$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
